$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 6 ("grandes regiões e unidades da federação"), shifting all rows below it up by one.
$ws.Rows.Item(6).Delete()
